# Case 4_1 (380 kV) line active-power-loss results: pl_mw.xlsx
# Updates the numeric results table in B2:O25 (columns F, J, K stay 0).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> column letter -> new value
$newValues = @{
    2 = @{ B=0.8167158455930235; C=0.243606436479098; D=0.07896030186675773; E=0.1168246504437014; G=0.4998711791772763; H=0.6319536109086101; I=0.5996809284984224; L=0.2021607255620097; M=0.1933652658690903; N=1.235775222517312; O=2.227114164706109 }
    3 = @{ B=0.735889592759861; C=0.231322646927282; D=0.07156802134579721; E=0.1178508831497402; G=0.496601782877562; H=0.6343799264883501; I=0.6049671442674232; L=0.199379554674735; M=0.180185750366455; N=1.245092473213312; O=2.224971769258616 }
    4 = @{ B=0.6863035602764853; C=0.2237174891204745; D=0.06706261964929183; E=0.1185188990756605; G=0.4949735912858415; H=0.6361843334324249; I=0.6085652076556833; L=0.1977714459579119; M=0.172143436586019; N=1.251271022627471; O=2.225108455026771 }
    5 = @{ B=0.6661086858632643; C=0.2206026875447549; D=0.06523506655933886; E=0.1188006642823547; G=0.4944053419961705; H=0.6369988018091135; I=0.6101199930828365; L=0.1971412247346862; M=0.1688789161764319; N=1.253904044391575; O=2.225529148757857 }
    6 = @{ B=0.6627560976664313; C=0.220084537596648; D=0.06493211259977727; E=0.1188480279186628; G=0.4943167344351806; H=0.6371388257661863; I=0.6103835106339162; L=0.1970380946578416; M=0.1683376240071439; N=1.254348217674007; O=2.225621044784049 }
    7 = @{ B=0.6860311553922713; C=0.2236755449100229; D=0.06703793845920814; E=0.1185226604034528; G=0.4949655421504673; H=0.6361949970793859; I=0.6085858175784793; L=0.1977628448734947; M=0.172099358047447; N=1.251306065832978; O=2.225112651014967 }
    8 = @{ B=0.788839114474257; C=0.2393841306022466; D=0.0764044990751529; E=0.1171706354877067; G=0.4986651140945355; H=0.6327249237267125; I=0.6014304676505766; L=0.2011811560239138; M=0.1888107295724311; N=1.238892879177307; O=2.226073977253606 }
    9 = @{ B=0.9907249798126827; C=0.2696837569978641; D=0.09503837246479918; E=0.1148195229011626; G=0.5089350954133351; H=0.6284153845717384; I=0.5901964181378361; L=0.208671955973962; M=0.221969955870847; N=1.218177531808209; O=2.239491735291779 }
    10 = @{ B=1.139167875544331; C=0.2916309076441905; D=0.1088927987650834; E=0.1132743752399614; G=0.5183287355613544; H=0.6267693479512388; I=0.5836515004502765; L=0.2146532581521541; M=0.2465603121515159; N=1.205162584370143; O=2.256400294778388 }
    11 = @{ B=1.206713474111325; C=0.3015458566225675; D=0.1152316541627272; E=0.1126108431791146; G=0.5230057625757496; H=0.6263504587191022; I=0.5810458089727284; L=0.2174776295043017; M=0.2577950332424805; N=1.199719189377767; O=2.265627992271476 }
    12 = @{ B=1.232292645149698; C=0.30529033220202; D=0.1176372535958308; E=0.1123652291364738; G=0.5248350453028223; H=0.6262392519639803; I=0.5801125989595874; L=0.2185619667716736; M=0.2620561015661877; N=1.197726443737601; O=2.269343411427542 }
    13 = @{ B=1.226783687031798; C=0.3044843446273831; D=0.117118933264635; E=0.1124178753256491; G=0.5244384863446498; H=0.6262610937705375; I=0.580311201509307; L=0.2183277777589723; M=0.2611381091113358; N=1.198152569409331; O=2.268533392777755 }
    14 = @{ B=1.208817874045849; C=0.3018541203960297; D=0.1154294599443375; E=0.112590523159386; G=0.5231550917521872; H=0.6263403594804657; I=0.5809679604294544; L=0.2175665421680293; M=0.2581454605452294; N=1.199553871519853; O=2.265929229764168 }
    15 = @{ B=1.197813396181232; C=0.3002417111061959; D=0.1143952868751086; E=0.1126970106153893; G=0.5223765578185464; H=0.6263950865256334; I=0.5813772147184402; L=0.2171021900003041; M=0.2563132452454653; N=1.200421134836589; O=2.264362903943208 }
    16 = @{ B=1.134753849781475; C=0.2909815393342114; D=0.1084792703333335; E=0.1133185299195452; G=0.5180312174269091; H=0.6268033599325662; I=0.5838292725118599; L=0.2144707547855376; M=0.2458270507612568; N=1.205527919135832; O=2.255828169074988 }
    17 = @{ B=1.096072452397095; C=0.2852829367180618; D=0.1048592946286391; E=0.1137098871782576; G=0.515469008380407; H=0.627138302479409; I=0.5854287621407934; L=0.2128829025575527; M=0.2394063449924317; N=1.208782931974966; O=2.250985932500527 }
    18 = @{ B=1.073825736246704; C=0.2819987754809574; D=0.1027806086105585; E=0.1139386911382076; G=0.5140332948411697; H=0.6273620057303333; I=0.5863837211721474; L=0.211979352853902; M=0.2357178992172493; N=1.210700046481826; O=2.248345359658714 }
    19 = @{ B=1.066293735256295; C=0.2808857084874035; D=0.1020773907309689; E=0.1140167969535817; G=0.5135537098000498; H=0.6274430818841807; I=0.5867130587380736; L=0.2116751014798552; M=0.2344698479118961; N=1.211356865180761; O=2.247476127265656 }
    20 = @{ B=1.100189977671846; C=0.2858902342844658; D=0.1052442925731896; E=0.1136678430800204; G=0.5157378259920335; H=0.6270994335451405; I=0.5852548734619631; L=0.2130509244719434; M=0.2400893684139547; N=1.20843178203161; O=2.251486434429864 }
    21 = @{ B=1.214094850363722; C=0.3026269566442181; D=0.1159255579342329; E=0.1125396590085543; G=0.5235304754839376; H=0.6263157904969603; I=0.5807736016474863; L=0.2177897341279476; M=0.2590242934404117; N=1.199140415208468; O=2.266688133433064 }
    22 = @{ B=1.288544097649265; C=0.3135064113266708; D=0.1229367532744163; E=0.1118352604641064; G=0.5289626796728868; H=0.6260800098028483; I=0.5781567537837944; L=0.2209731082149915; M=0.2714384372108967; N=1.193467503022305; O=2.277912110096366 }
    23 = @{ B=1.248809093646969; C=0.3077052991997107; D=0.1191919752222503; E=0.1122082007147664; G=0.5260323266756188; H=0.626180569271753; I=0.5795248508941135; L=0.2192662081753838; M=0.2648092798858741; N=1.196458706703432; O=2.271803663741252 }
    24 = @{ B=1.098328470126887; C=0.2856156998319932; D=0.105070227257599; E=0.1136868393564754; G=0.5156161772651586; H=0.6271169091876914; I=0.585333378297392; L=0.2129749327051087; M=0.239780564462734; N=1.208590394361956; O=2.251259711173958 }
    25 = @{ B=0.9360850605334576; C=0.2615415967271417; D=0.08996872981231263; E=0.1154235049371408; G=0.505832945658355; H=0.6293141973740433; I=0.5929356889483515; L=0.2065614229862476; M=0.2129588110820109; N=1.223393979742141; O=2.234625159975309 }
}

foreach ($row in $newValues.Keys) {
    $rowData = $newValues[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}

Write-Host "Updated pl_mw.xlsx results (rows 2-25, cols B-O) for the 380 kV case"
